$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.363.20"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.563.72"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.06"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3736"
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3263"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.31"
$ws.Range("E9").Value = "  -5.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.140"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07380"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.35"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.816"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.801"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.558.20"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001093"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06714"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.53"
$ws.Range("E19").Value = "  -2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.327"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.19"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.63"
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.357.75"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.520"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.67"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.34"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.905"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.81"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.736.94"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.049"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.936"
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.882"
$ws.Range("E34").Value = "  -4.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.508"
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08200"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02378"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06283"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.280"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2175"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.234"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.96"
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6059"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.57"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.735"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5880"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.987"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.16"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.174"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07135"
